{"js": "// \"Add Figure (from Pandoc 3)\" \u2014 bring the Pandoc reference-docx style\n// definitions up to date. The concrete, semantically meaningful deltas\n// (ignoring pure XML-attribute/element reordering noise from the\n// canonicalizer) are:\n//\n//   1. Add a new paragraph style \"Abstract Title\" (styleId AbstractTitle),\n//      based on Normal, followed by Abstract.\n//   2. Tighten the existing \"Abstract\" style's space-before from 300 to\n//      100 twips (15pt -> 5pt); space-after stays 300 (15pt).\n//   3. Add a new paragraph style \"Footnote Block Text\"\n//      (styleId FootnoteBlockText), based on \"Footnote Text\".\n//   4. Give the \"ImportTok\" character style bold + green (#008000) text.\n//   5. Give the \"BuiltInTok\" character style green (#008000) text.\n\nconst styles = context.document.getStyles();\n\n// --- 1. New \"Abstract Title\" paragraph style --------------------------\ncontext.document.addStyle(\"Abstract Title\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst abstractTitle = styles.getByName(\"Abstract Title\");\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\n\nconst atPf = abstractTitle.paragraphFormat;\natPf.keepWithNext = true;\natPf.keepTogether = true;\natPf.spaceBefore = 15; // pt  (300 twips)\natPf.spaceAfter = 0;\natPf.alignment = Word.Alignment.centered;\n\nconst atFont = abstractTitle.font;\natFont.bold = true;\natFont.color = \"#345A8A\";\natFont.size = 10; // pt (sz 20 half-points)\natFont.sizeBidirectional = 10; // szCs 20\n\nawait context.sync();\n\n// --- 2. \"Abstract\" style: space-before 300 -> 100 twips (15pt -> 5pt) --\nconst abstract = styles.getByName(\"Abstract\");\nabstract.paragraphFormat.spaceBefore = 5; // pt (100 twips)\nawait context.sync();\n\n// --- 3. New \"Footnote Block Text\" paragraph style ----------------------\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst fnBlockText = styles.getByName(\"Footnote Block Text\");\nfnBlockText.baseStyle = \"Footnote Text\";\nfnBlockText.nextParagraphStyle = \"Footnote Text\";\nfnBlockText.priority = 9;\nfnBlockText.unhideWhenUsed = true;\nfnBlockText.quickStyle = true;\n\nconst fnPf = fnBlockText.paragraphFormat;\nfnPf.spaceBefore = 5; // pt (100 twips)\nfnPf.spaceAfter = 5; // pt (100 twips)\nfnPf.firstLineIndent = 0;\nfnPf.leftIndent = 24; // pt (480 twips)\nfnPf.rightIndent = 24; // pt (480 twips)\n\nawait context.sync();\n\n// --- 4. \"ImportTok\" character style: bold + green text -----------------\nconst importTok = styles.getByName(\"ImportTok\");\nimportTok.font.bold = true;\nimportTok.font.color = \"#008000\";\nawait context.sync();\n\n// --- 5. \"BuiltInTok\" character style: green text ------------------------\nconst builtInTok = styles.getByName(\"BuiltInTok\");\nbuiltInTok.font.color = \"#008000\";\nawait context.sync();\n", "ps1": "# \"Add Figure (from Pandoc 3)\" \u2014 bring the Pandoc reference-docx style\n# definitions up to date. The concrete, semantically meaningful deltas\n# (ignoring pure XML-attribute/element reordering noise from the\n# canonicalizer) are:\n#\n#   1. Add a new paragraph style \"Abstract Title\" (styleId AbstractTitle),\n#      based on Normal, followed by Abstract.\n#   2. Tighten the existing \"Abstract\" style's space-before from 300 to\n#      100 twips (15pt -> 5pt); space-after stays 300 (15pt).\n#   3. Add a new paragraph style \"Footnote Block Text\"\n#      (styleId FootnoteBlockText), based on \"Footnote Text\".\n#   4. Give the \"ImportTok\" character style bold + green (#008000) text.\n#   5. Give the \"BuiltInTok\" character style green (#008000) text.\n\n$d = $word.ActiveDocument\n\n$wdStyleTypeParagraph = 1\n$wdAlignParagraphCenter = 1\n\n# --- 1. New \"Abstract Title\" paragraph style ---------------------------\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", $wdStyleTypeParagraph)\n$abstractTitle.BaseStyle = \"Normal\"\n$abstractTitle.NextParagraphStyle = \"Abstract\"\n$abstractTitle.QuickStyle = $true\n\n$atPf = $abstractTitle.ParagraphFormat\n$atPf.KeepWithNext = $true\n$atPf.KeepTogether = $true\n$atPf.SpaceBefore = 15   # pt (300 twips)\n$atPf.SpaceAfter = 0\n$atPf.Alignment = $wdAlignParagraphCenter\n\n$atFont = $abstractTitle.Font\n$atFont.Bold = $true\n$atFont.Color = 9067060  # BGR for RGB 345A8A\n$atFont.Size = 10        # pt (sz 20 half-points)\n$atFont.SizeBi = 10      # szCs 20\n\n# --- 2. \"Abstract\" style: space-before 300 -> 100 twips (15pt -> 5pt) --\n$abstract = $d.Styles(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5   # pt (100 twips)\n\n# --- 3. New \"Footnote Block Text\" paragraph style -----------------------\n$fnBlockText = $d.Styles.Add(\"Footnote Block Text\", $wdStyleTypeParagraph)\n$fnBlockText.BaseStyle = \"Footnote Text\"\n$fnBlockText.NextParagraphStyle = \"Footnote Text\"\n$fnBlockText.Priority = 9\n$fnBlockText.UnhideWhenUsed = $true\n$fnBlockText.QuickStyle = $true\n\n$fnPf = $fnBlockText.ParagraphFormat\n$fnPf.SpaceBefore = 5     # pt (100 twips)\n$fnPf.SpaceAfter = 5      # pt (100 twips)\n$fnPf.FirstLineIndent = 0\n$fnPf.LeftIndent = 24     # pt (480 twips)\n$fnPf.RightIndent = 24    # pt (480 twips)\n\n# --- 4. \"ImportTok\" character style: bold + green text ------------------\n$importTok = $d.Styles(\"ImportTok\")\n$importTok.Font.Bold = $true\n$importTok.Font.Color = 32768  # BGR for RGB 008000\n\n# --- 5. \"BuiltInTok\" character style: green text -------------------------\n$builtInTok = $d.Styles(\"BuiltInTok\")\n$builtInTok.Font.Color = 32768  # BGR for RGB 008000\n"}
